# "Added more stories to documentation"
#
# The story-breakdown bullet list currently ends with one empty list
# item. We turn that empty item into a new story, and append one more
# story bullet after it. The "_GoBack" bookmark (which Word drops at the
# location of the last edit) needs to move from the old last-edited
# paragraph to the end of the text we just typed.

$d = $word.ActiveDocument

# 1) Detach the stale _GoBack bookmark from its old position (end of the
#    "handle the master/storefront catalog relationship" bullet).
$goBack = $d.Bookmarks("_GoBack")
$goBack.Delete()

# 2) The bullet list currently ends with one empty "ListParagraph" item;
#    give it the first new story's text.
$lastPara = $d.Paragraphs.Last
$lastRange = $lastPara.Range
$lastRange.InsertAfter("I want to have ability to submit comma delimited list of products to export")

# 3) Add a new list paragraph after it (inherits the ListParagraph /
#    numbering formatting) for the second new story. Type a temporary
#    marker after the text so the insertion point used for the bookmark
#    below is not the very last character in the document, then strip
#    the marker back out once the bookmark is anchored.
$lastRange.InsertParagraphAfter()
$newPara = $d.Paragraphs.Last
$newRange = $newPara.Range
$newRange.InsertAfter("I want to click link to easily download exported catalogZZMARKERZZ")

$markerRange = $d.Content
$markerRange.Find.Execute("ZZMARKERZZ") | Out-Null
$bookmarkSpot = $d.Range($markerRange.Start, $markerRange.Start)
$d.Bookmarks.Add("_GoBack", $bookmarkSpot)
$markerRange.Text = ""
